# Update LR-pair TPM-derived values for Bmp15-Acvr2a (OldD7) per new TPM script run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster = ECs)
$ws.Range("G2").Value = 0.053572
$ws.Range("H2").Value = 0.160716
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 0.8726047005413334
$ws.Range("R2").Value = 7.853442304872
$ws.Range("S2").Value = 0.2176904746803693
$ws.Range("T2").Value = 0.2176904746803693

# Row 3 (Target cluster = FAPs)
$ws.Range("G3").Value = 0.053572
$ws.Range("H3").Value = 0.160716
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83271999999999
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 1.479171491946667
$ws.Range("R3").Value = 13.31254342752
$ws.Range("S3").Value = 0.3690119294748028
$ws.Range("T3").Value = 0.3690119294748029

# Row 4 (Target cluster = MuSCs)
$ws.Range("G4").Value = 0.053572
$ws.Range("H4").Value = 0.160716
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 1.4071395629
$ws.Range("R4").Value = 12.6642560661
$ws.Range("S4").Value = 0.3510419771967738
$ws.Range("T4").Value = 0.3510419771967739

# Row 5 (Target cluster = Resolving-Mac)
$ws.Range("G5").Value = 0.053572
$ws.Range("H5").Value = 0.160716
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 0.2495494832613333
$ws.Range("R5").Value = 2.245945349352
$ws.Range("S5").Value = 0.06225561864805391
$ws.Range("T5").Value = 0.06225561864805392
